$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $TextValue)
    $c = $ws.Range($CellRef)
    $c.NumberFormat = "@"
    $c.Value = $TextValue
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "69.174.47"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "3.924.94"
$ws.Range("E3").Value = "  +5.01%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.15%  "
Set-TextValue "D5" "606.38"
$ws.Range("E5").Value = "  +0.86%  "
Set-TextValue "D6" "165.05"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "3.920.90"
$ws.Range("E7").Value = "  +4.98%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  -1.40%  "
Set-TextValue "D11" "6.41"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  +0.98%  "
Set-TextValue "D13" "37.25"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "4.579.11"
$ws.Range("E15").Value = "  +4.91%  "
$ws.Range("D16").Value = "3.920.36"
$ws.Range("E16").Value = "  +5.18%  "
$ws.Range("D17").Value = "69.215.29"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("E19").Value = "  -0.28%  "
Set-TextValue "D20" "17.19"
$ws.Range("E20").Value = "  -2.66%  "
Set-TextValue "D21" "11.22"
$ws.Range("E21").Value = "  +0.37%  "
Set-TextValue "D22" "489.30"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +12.86%  "
Set-TextValue "D25" "84.53"
$ws.Range("E25").Value = "  -0.09%  "
Set-TextValue "D26" "2.27"
$ws.Range("E26").Value = "  -0.15%  "
Set-TextValue "D27" "12.18"
$ws.Range("E27").Value = "  -1.32%  "
Set-TextValue "D28" "10.12"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").Value = "4.075.72"
$ws.Range("E31").Value = "  +4.94%  "
Set-TextValue "D32" "7.91"
$ws.Range("E32").Value = "  -3.29%  "
Set-TextValue "D33" "32.46"
$ws.Range("E33").Value = "  +2.96%  "
Set-TextValue "D34" "2.39"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").Value = "3.868.81"
$ws.Range("E35").Value = "  +5.31%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("E37").Value = "  +3.34%  "
$ws.Range("E38").Value = "  +1.94%  "
Set-TextValue "D39" "5.96"
$ws.Range("E39").Value = "  +0.62%  "
Set-TextValue "D40" "0.998"
$ws.Range("E40").Value = "  -0.18%  "
Set-TextValue "D41" "0.322"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("E42").Value = "  -3.55%  "
Set-TextValue "D43" "440.85"
$ws.Range("E43").Value = "  +3.47%  "
Set-TextValue "D44" "2.01"
$ws.Range("E44").Value = "  +1.21%  "
Set-TextValue "D45" "48.49"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D46" "8.48"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D47" "1.00"
$ws.Range("E47").Value = "  +0.01%  "
Set-TextValue "D48" "27.80"
$ws.Range("E48").Value = "  +18.34%  "
$ws.Range("D49").Value = "2.849.49"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D50" "141.56"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D51" "0.0359"
$ws.Range("E51").Value = "  +2.34%  "
